# Update gh-pages to output generated at 456a3b4
# Refreshes "想去人数" (F) / "最低票价" (G) figures on the 展览 (Exhibition)
# sheet and mirrors the same refresh onto 全部类型 (All types), which is a
# concatenation of 展览 followed by 演出.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows 2-24) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 1596
$ws1.Range("F3").Value = 8956
$ws1.Range("G4").Value = 89
$ws1.Range("F5").Value = 500
$ws1.Range("F6").Value = 682
$ws1.Range("F7").Value = 343
$ws1.Range("F9").Value = 44
$ws1.Range("F10").Value = 71
$ws1.Range("F11").Value = 3816
$ws1.Range("F14").Value = 100
$ws1.Range("F15").Value = 4208
$ws1.Range("F16").Value = 6
$ws1.Range("F18").Value = 1140
$ws1.Range("F19").Value = 6
$ws1.Range("F21").Value = 239
$ws1.Range("F22").Value = 10
$ws1.Range("F23").Value = 2618
$ws1.Range("F24").Value = 105

# --- Sheet "全部类型" (展览 rows 2-24 map to rows 2-24, then 演出 row 2 maps
#     to row 25, same F/G edits as 展览 except F15's prior value differed
#     and the last 展览 row lands on row 25 here) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 1596
$ws4.Range("F3").Value = 8956
$ws4.Range("G4").Value = 89
$ws4.Range("F5").Value = 500
$ws4.Range("F6").Value = 682
$ws4.Range("F7").Value = 343
$ws4.Range("F9").Value = 44
$ws4.Range("F10").Value = 71
$ws4.Range("F11").Value = 3816
$ws4.Range("F14").Value = 100
$ws4.Range("F15").Value = 4208
$ws4.Range("F16").Value = 6
$ws4.Range("F18").Value = 1140
$ws4.Range("F19").Value = 6
$ws4.Range("F21").Value = 239
$ws4.Range("F22").Value = 10
$ws4.Range("F23").Value = 2618
$ws4.Range("F25").Value = 105
